# Add a new "Fullname" column to the category template sheet.
# This inserts a column before the existing "Name" column (B) and
# populates its header + sample row, shifting Name/Img/Desc/Era/Franchise
# one column to the right (B->C, C->D, D->E, E->F, F->G) and extending the
# Row 1 title merge from A1:E1 to A1:F1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing Name/Img/Desc/Era/Franchise columns one to the right by
# inserting a brand new column at B (the mergeCell A1:E1 title bar grows
# along with it automatically).
$ws.Columns("B").Insert()

# New "Fullname" column: header + sample value copied from the Name column
# (now in C after the shift).
$ws.Range("B2").Value = "Fullname"
$ws.Range("B3").Value = $ws.Range("C3").Value2

# Match the new column's width to the authored template.
$ws.Columns("B").ColumnWidth = 13.8333333333

# Keep selection on the new cell, matching the authored workbook.
$ws.Range("B3").Select()
